$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.423.31"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.565.70"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'287.87"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("D7").Value = "'0.3725"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D8").Value = "'48.19"
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("D9").Value = "'0.3310"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.129"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.07459"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'20.57"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").Value = "'5.919"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "'6.895"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "1.562.12"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'0.00001110"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "'0.06756"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'87.62"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D21").Value = "'6.333"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "'16.36"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "'12.05"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "22.413.41"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'2.383"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").Value = "'2.555"
$ws.Range("E26").Value = "  -4.23%  "
$ws.Range("D27").Value = "'153.09"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").Value = "'19.61"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("D29").Value = "'5.012"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "'123.98"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "1.738.90"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").Value = "'1.052"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'6.103"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").Value = "'9.611"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").Value = "'0.08303"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").Value = "'0.02449"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("D38").Value = "'0.2267"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").Value = "'0.06370"
$ws.Range("E39").Value = "  -2.63%  "
$ws.Range("D40").Value = "'1.283"
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("D41").Value = "'5.315"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").Value = "'0.6264"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").Value = "'1.003"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'13.73"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").Value = "'0.6107"
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("D47").Value = "'3.767"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").Value = "'2.036"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").Value = "'125.12"
$ws.Range("D50").Value = "'1.208"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("D51").Value = "'0.07221"
$ws.Range("E51").Value = "  -1.11%  "
